$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.174.35'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.36%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.734.51'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.33%  '

$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '622.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.60%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.99%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.728.04'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.43%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('E9').Value = '  -1.38%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.168'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.08%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.31'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.58%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.486'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.26%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.72'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.43%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000260'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.11%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.359.16'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.23%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.734.93'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.04%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.227.38'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.36%  '

$ws.Range('E18').Value = '  -1.57%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.60%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.90%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '506.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.42%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.29'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.725'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.04%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.56'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.62%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.74'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.20%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.48'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.43%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '13.14'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.20%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000137'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +21.04%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.13%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.49'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.17%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.94'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.60%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.96'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.68%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.18'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.16%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.115'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.85%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.06%  '

$ws.Range('E36').Value = '  +0.88%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.17'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.46%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.138'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.50%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.339'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.79%  '

$ws.Range('E40').Value = '  -6.26%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.47'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.47%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '45.71'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.22%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '435.05'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.45%  '

$ws.Range('E44').Value = '  +1.42%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.71'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.73%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.012.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.93%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0364'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.43%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.28%  '

$ws.Range('E49').Value = '  -0.03%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '137.90'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.14%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.50'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.40%  '
